# Generate Report for Handoff
# Updates the localization-status workbook: the fff5f01f-a8a7-45ec-b9dd-12e41eac4dde
# source file has moved from "In Translation" to "Ready for handoff", with a refreshed
# handoff timestamp and priority, reflected across the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-09-06 00:15:37"

# Widen the zh-cn / de-de status columns to fit the new, longer status text.
$overview.Columns.Item(5).ColumnWidth = 16.333333333333336
$overview.Columns.Item(6).ColumnWidth = 16.333333333333336

# --- zh-cn sheet ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("E3").Value = "mt"
$zhcn.Range("H3").Value = "2016-09-06 00:15:32"
$zhcn.Columns.Item(3).ColumnWidth = 16.333333333333336

# --- de-de sheet ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("E3").Value = "mt"
$dede.Range("H3").Value = "2016-09-06 00:15:37"
$dede.Columns.Item(3).ColumnWidth = 16.333333333333336
